$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 8
$ws.Range("H8").Value = 60
$ws.Range("I8").Value = 60
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0

# Row 10
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0

# Row 11
$ws.Range("H11").Value = 80
$ws.Range("I11").Value = 45
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0

# Row 12
$ws.Range("H12").Value = 80
$ws.Range("I12").Value = 45
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0

# Row 13
$ws.Range("H13").Value = 60
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0

# Row 14
$ws.Range("O14").Value = 90
$ws.Range("P14").Value = 90
$ws.Range("Q14").Value = 90
$ws.Range("R14").Value = 45
$ws.Range("S14").Value = 0

# Row 15
$ws.Range("O15").Value = 120
$ws.Range("P15").Value = 120
$ws.Range("Q15").Value = 120
$ws.Range("R15").Value = 30
$ws.Range("S15").Value = 0

# Row 16
$ws.Range("H16").Value = 60
$ws.Range("I16").Value = 60
$ws.Range("J16").Value = 60
$ws.Range("K16").Value = 0

# Update sheet view - topLeftCell and active cell selection
$ws.Range("L20").Select()

$wb.Save()
